# Rearrange the 9 stacked chart images on the "Graphs" sheet into a 3x3 grid.
# Before: all 9 pictures anchored at column 0, rows 0,25,50,75,100,125,150,175,200
# After : pictures anchored on a 3-column (0,10,20) x 3-row (0,20,40) grid:
#
#           col 0      col 10     col 20
# row  0 :  Image 1    Image 2    Image 3
# row 20 :  Image 4    Image 5    Image 6
# row 40 :  Image 7    Image 8    Image 9
#
# Shapes stay anchored to their top-left cell (oneCellAnchor) when moved via
# Left/Top, so we convert the target (col, row) grid cell into point offsets
# using the sheet's (uniform) column width / row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colWidthPts = $ws.Columns.Item(1).Width
$rowHeightPts = $ws.Rows.Item(1).Height

# Shape index (1-based, matches "Image N" creation order) -> (target col, target row), both 0-based.
$targets = @{
    1 = @(0, 0)
    2 = @(10, 0)
    3 = @(20, 0)
    4 = @(0, 20)
    5 = @(10, 20)
    6 = @(20, 20)
    7 = @(0, 40)
    8 = @(10, 40)
    9 = @(20, 40)
}

foreach ($idx in $targets.Keys) {
    $target = $targets[$idx]
    $targetCol = $target[0]
    $targetRow = $target[1]

    $shp = $ws.Shapes.Item([int]$idx)
    $shp.Left = $targetCol * $colWidthPts
    $shp.Top = $targetRow * $rowHeightPts
}

# The new layout reaches down to row 41 / column U (1-based), so the rows that
# anchor a picture (1, 21, 41) get touched so they're materialized in the
# sheet's row list.
$ws.Rows.Item(1).OutlineLevel = 0
$ws.Rows.Item(21).OutlineLevel = 0
$ws.Rows.Item(41).OutlineLevel = 0
